$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a "last changed" date serial. The whole
# column was refreshed by one day (45207 -> 45208) for every data row.
$lastRow = $ws.UsedRange.Rows.Count
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value = 45208
    }
}
